# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps on the zh-cn and de-de report sheets, as part of regenerating
# the handback status report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-24 20:26:12"
$wsZhCn.Range("H2").Value = "2016-03-24 20:26:44"
$wsZhCn.Range("E4").Value = "2016-03-24 20:26:12"
$wsZhCn.Range("H4").Value = "2016-03-24 20:26:44"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-24 20:26:17"
$wsDeDe.Range("H2").Value = "2016-03-24 20:26:51"
$wsDeDe.Range("E4").Value = "2016-03-24 20:26:17"
$wsDeDe.Range("H4").Value = "2016-03-24 20:26:51"
